$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.04537966666667
$ws.Range("H2").Value = 54.13613900000001
$ws.Range("I2").Value = 0.6797959733292525
$ws.Range("J2").Value = 0.6797959733292525
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.667552
$ws.Range("N2").Value = 263.002656
$ws.Range("O2").Value = 0.3606416352150456
$ws.Range("P2").Value = 0.3606416352150456
$ws.Range("Q2").Value = 1581.994260287243
$ws.Range("R2").Value = 14237.94834258519
$ws.Range("S2").Value = 0.2451627314340651
$ws.Range("T2").Value = 0.2451627314340651

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.04537966666667
$ws.Range("H3").Value = 54.13613900000001
$ws.Range("I3").Value = 0.6797959733292525
$ws.Range("J3").Value = 0.6797959733292525
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 72.97955566666666
$ws.Range("N3").Value = 218.938667
$ws.Range("O3").Value = 0.3002190170987564
$ws.Range("P3").Value = 0.3002190170987564
$ws.Range("Q3").Value = 1316.943789909635
$ws.Range("R3").Value = 11852.49410918671
$ws.Range("S3").Value = 0.2040876789406006
$ws.Range("T3").Value = 0.2040876789406006

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.04537966666667
$ws.Range("H4").Value = 54.13613900000001
$ws.Range("I4").Value = 0.6797959733292525
$ws.Range("J4").Value = 0.6797959733292525
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 41.37117366666666
$ws.Range("N4").Value = 124.113521
$ws.Range("O4").Value = 0.1701903085181653
$ws.Range("P4").Value = 0.1701903085181653
$ws.Range("Q4").Value = 746.5585360706021
$ws.Range("R4").Value = 6719.02682463542
$ws.Range("S4").Value = 0.1156946864303119
$ws.Range("T4").Value = 0.1156946864303119

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.04537966666667
$ws.Range("H5").Value = 54.13613900000001
$ws.Range("I5").Value = 0.6797959733292525
$ws.Range("J5").Value = 0.6797959733292525
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 41.06943633333334
$ws.Range("N5").Value = 123.208309
$ws.Range("O5").Value = 0.1689490391680327
$ws.Range("P5").Value = 0.1689490391680327
$ws.Range("Q5").Value = 741.1135713309947
$ws.Range("R5").Value = 6670.022141978952
$ws.Range("S5").Value = 0.1148508765242748
$ws.Range("T5").Value = 0.1148508765242748

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6001993333333334
$ws.Range("H6").Value = 1.800598
$ws.Range("I6").Value = 0.02261039099934159
$ws.Range("J6").Value = 0.02261039099934159
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 87.667552
$ws.Range("N6").Value = 263.002656
$ws.Range("O6").Value = 0.3606416352150456
$ws.Range("P6").Value = 0.3606416352150456
$ws.Range("Q6").Value = 52.61800626536534
$ws.Range("R6").Value = 473.5620563882881
$ws.Range("S6").Value = 0.008154248382854102
$ws.Range("T6").Value = 0.008154248382854102

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6001993333333334
$ws.Range("H7").Value = 1.800598
$ws.Range("I7").Value = 0.02261039099934159
$ws.Range("J7").Value = 0.02261039099934159
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 72.97955566666666
$ws.Range("N7").Value = 218.938667
$ws.Range("O7").Value = 0.3002190170987564
$ws.Range("P7").Value = 0.3002190170987564
$ws.Range("Q7").Value = 43.80228065809622
$ws.Range("R7").Value = 394.220525922866
$ws.Range("S7").Value = 0.006788069362040901
$ws.Range("T7").Value = 0.006788069362040901

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.6001993333333334
$ws.Range("H8").Value = 1.800598
$ws.Range("I8").Value = 0.02261039099934159
$ws.Range("J8").Value = 0.02261039099934159
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 41.37117366666666
$ws.Range("N8").Value = 124.113521
$ws.Range("O8").Value = 0.1701903085181653
$ws.Range("P8").Value = 0.1701903085181653
$ws.Range("Q8").Value = 24.83095085395089
$ws.Range("R8").Value = 223.478557685558
$ws.Range("S8").Value = 0.003848069419894294
$ws.Range("T8").Value = 0.003848069419894294

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.6001993333333334
$ws.Range("H9").Value = 1.800598
$ws.Range("I9").Value = 0.02261039099934159
$ws.Range("J9").Value = 0.02261039099934159
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 41.06943633333334
$ws.Range("N9").Value = 123.208309
$ws.Range("O9").Value = 0.1689490391680327
$ws.Range("P9").Value = 0.1689490391680327
$ws.Range("Q9").Value = 24.64984830764245
$ws.Range("R9").Value = 221.848634768782
$ws.Range("S9").Value = 0.003820003834552296
$ws.Range("T9").Value = 0.003820003834552297

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.690054333333333
$ws.Range("H10").Value = 14.070163
$ws.Range("I10").Value = 0.1766812397072912
$ws.Range("J10").Value = 0.1766812397072912
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 87.667552
$ws.Range("N10").Value = 263.002656
$ws.Range("O10").Value = 0.3606416352150456
$ws.Range("P10").Value = 0.3606416352150456
$ws.Range("Q10").Value = 411.1655821503253
$ws.Range("R10").Value = 3700.490239352928
$ws.Range("S10").Value = 0.06371861119985893
$ws.Range("T10").Value = 0.06371861119985894

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.690054333333333
$ws.Range("H11").Value = 14.070163
$ws.Range("I11").Value = 0.1766812397072912
$ws.Range("J11").Value = 0.1766812397072912
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 72.97955566666666
$ws.Range("N11").Value = 218.938667
$ws.Range("O11").Value = 0.3002190170987564
$ws.Range("P11").Value = 0.3002190170987564
$ws.Range("Q11").Value = 342.2780812991912
$ws.Range("R11").Value = 3080.502731692721
$ws.Range("S11").Value = 0.05304306812471272
$ws.Range("T11").Value = 0.05304306812471273

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.690054333333333
$ws.Range("H12").Value = 14.070163
$ws.Range("I12").Value = 0.1766812397072912
$ws.Range("J12").Value = 0.1766812397072912
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 41.37117366666666
$ws.Range("N12").Value = 124.113521
$ws.Range("O12").Value = 0.1701903085181653
$ws.Range("P12").Value = 0.1701903085181653
$ws.Range("Q12").Value = 194.0330523304359
$ws.Range("R12").Value = 1746.297470973923
$ws.Range("S12").Value = 0.0300694346951558
$ws.Range("T12").Value = 0.0300694346951558

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.690054333333333
$ws.Range("H13").Value = 14.070163
$ws.Range("I13").Value = 0.1766812397072912
$ws.Range("J13").Value = 0.1766812397072912
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 41.06943633333334
$ws.Range("N13").Value = 123.208309
$ws.Range("O13").Value = 0.1689490391680327
$ws.Range("P13").Value = 0.1689490391680327
$ws.Range("Q13").Value = 192.6178878427075
$ws.Range("R13").Value = 1733.560990584367
$ws.Range("S13").Value = 0.0298501256875637
$ws.Range("T13").Value = 0.02985012568756371

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.209654333333333
$ws.Range("H14").Value = 9.628962999999999
$ws.Range("I14").Value = 0.1209123959641148
$ws.Range("J14").Value = 0.1209123959641148
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 87.667552
$ws.Range("N14").Value = 263.002656
$ws.Range("O14").Value = 0.3606416352150456
$ws.Range("P14").Value = 0.3606416352150456
$ws.Range("Q14").Value = 281.3825381695253
$ws.Range("R14").Value = 2532.442843525728
$ws.Range("S14").Value = 0.04360604419826744
$ws.Range("T14").Value = 0.04360604419826744

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.209654333333333
$ws.Range("H15").Value = 9.628962999999999
$ws.Range("I15").Value = 0.1209123959641148
$ws.Range("J15").Value = 0.1209123959641148
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 72.97955566666666
$ws.Range("N15").Value = 218.938667
$ws.Range("O15").Value = 0.3002190170987564
$ws.Range("P15").Value = 0.3002190170987564
$ws.Range("Q15").Value = 234.2391470902578
$ws.Range("R15").Value = 2108.15232381232
$ws.Range("S15").Value = 0.03630020067140218
$ws.Range("T15").Value = 0.03630020067140218

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.209654333333333
$ws.Range("H16").Value = 9.628962999999999
$ws.Range("I16").Value = 0.1209123959641148
$ws.Range("J16").Value = 0.1209123959641148
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 41.37117366666666
$ws.Range("N16").Value = 124.113521
$ws.Range("O16").Value = 0.1701903085181653
$ws.Range("P16").Value = 0.1701903085181653
$ws.Range("Q16").Value = 132.7871668343025
$ws.Range("R16").Value = 1195.084501508723
$ws.Range("S16").Value = 0.02057811797280326
$ws.Range("T16").Value = 0.02057811797280326

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.209654333333333
$ws.Range("H17").Value = 9.628962999999999
$ws.Range("I17").Value = 0.1209123959641148
$ws.Range("J17").Value = 0.1209123959641148
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 41.06943633333334
$ws.Range("N17").Value = 123.208309
$ws.Range("O17").Value = 0.1689490391680327
$ws.Range("P17").Value = 0.1689490391680327
$ws.Range("Q17").Value = 131.8186942948408
$ws.Range("R17").Value = 1186.368248653567
$ws.Range("S17").Value = 0.02042803312164191
$ws.Range("T17").Value = 0.02042803312164191
